$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Nuevo articulo: id 999, "arma calibre 38", precio anterior 50000, precio nuevo 50000
$ws.Range("A17").Value = 999
$ws.Range("B17").Value = "arma calibre 38"
$ws.Range("C17").Value = 50000
$ws.Range("D17").Value = 50000

# Match the new cell cursor/selection left in the sheet after the edit
[void]$ws.Range("D18").Select()
